$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells F1 (height) and G1 (weight), matching style of existing headers (B1:E1)
$ws.Range("F1").Value = "height"
$ws.Range("G1").Value = "weight"

# Copy the formatting (style) from E1 onto F1:G1
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in data for rows 2-9
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.416666666666667   # column E: fantasy points
    $ws.Cells.Item($r, 6).Value = 255                  # column F: height
    $ws.Cells.Item($r, 7).Value = 0                    # column G: weight
}
